$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.851.97"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.405.13"
$ws.Range("E3").Value = "  +4.87%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'336.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.42%  "
$ws.Range("D6").Value = "'101.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.55%  "
$ws.Range("D7").Value = "'0.645"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.74%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.640"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.79%  "
$ws.Range("D10").Value = "'40.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.96%  "
$ws.Range("D11").Value = "'0.0931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "'8.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.36%  "
$ws.Range("E13").Value = "  -3.15%  "
$ws.Range("D14").Value = "'16.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.45%  "
$ws.Range("D15").Value = "'0.107"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").Value = "2.768.52"
$ws.Range("E16").Value = "  +4.95%  "
$ws.Range("D17").Value = "2.390.62"
$ws.Range("E17").Value = "  +4.35%  "
$ws.Range("D18").Value = "42.869.72"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'7.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.80%  "
$ws.Range("D20").Value = "'0.0000108"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "'3.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.07%  "
$ws.Range("D22").Value = "'76.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").Value = "'270.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.58%  "
$ws.Range("E24").Value = "  -4.69%  "
$ws.Range("D25").Value = "'10.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +15.31%  "
$ws.Range("D26").Value = "'11.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D28").Value = "'24.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.18%  "
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("D30").Value = "'174.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("D31").Value = "'3.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("D32").Value = "'0.0924"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("D33").Value = "'36.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.60%  "
$ws.Range("D34").Value = "'6.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.21%  "
$ws.Range("D35").Value = "'0.135"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.71%  "
$ws.Range("E36").Value = "  -7.00%  "
$ws.Range("D37").Value = "'0.0363"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.80%  "
$ws.Range("D38").Value = "'3.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.53%  "
$ws.Range("E39").Value = "  +2.58%  "
$ws.Range("D40").Value = "'2.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.80%  "
$ws.Range("D41").Value = "'1.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.38%  "
$ws.Range("D42").Value = "'0.234"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").Value = "'69.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "'93.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +46.29%  "
$ws.Range("D46").Value = "'118.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.77%  "
$ws.Range("D47").Value = "'11.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("D48").Value = "'5.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.30%  "
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "1.643.17"
$ws.Range("E50").Value = "  +11.05%  "
$ws.Range("D51").Value = "'1.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.90%  "
